$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.3908593446462402

$ws.Range("A3").Value = 20
$ws.Range("B3").Value = 0.4771933468515769

$ws.Range("A4").Value = 30
$ws.Range("B4").Value = 0.5171777931892586

$ws.Range("A5").Value = 40
$ws.Range("B5").Value = 0.5429371701380643
